$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (column I) and DialogAct (column J) values for the affected rows.

$updates = @(
    @{ Row = 18;  I = "%";  J = "Uninterpretable" },
    @{ Row = 48;  I = "ba"; J = "Appreciation" },
    @{ Row = 50;  I = "ba"; J = "Appreciation" },
    @{ Row = 54;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 67;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 72;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 86;  I = "%";  J = "Uninterpretable" },
    @{ Row = 89;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 101; I = "ba"; J = "Appreciation" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
